# Update the table of ERGM terms
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix wording: "same ritual treatment" -> "the same ritual treatment"
$ws.Range("C5").Value = "Density of ties between nodes with the same ritual treatment"

# Fix wording: "each pairs of burials" -> "each pair of burials"
$ws.Range("C9").Value = "Distance (in meter) between each pair of burials"

# Update the active cell selection to C5
$ws.Range("C5").Select()
